# Edit script: applies the changes described by the diff to before.pptx
# 1) Update cached datetimeFigureOut text ("10/29/23" -> "10/16/24") on the
#    slide master and every slide layout's Date placeholder.
# 2) Rewrite the "Reminder:" bullets on slide 1 into the new
#    "There’s a conference!:" bullets (merging/removing the old 3rd bullet).

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -ne $newText) {
                $tr.Text = $newText
            }
        }
    }
}

# --- 1) Slide master date placeholder ---
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes "10/16/24"

# --- 1b) Every custom (slide) layout's date placeholder ---
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes "10/16/24"
}

# --- 2) Slide 1 body text updates ---
$slide1 = $p.Slides.Item(1)
$body = $slide1.Shapes.Item(2)
$tr = $body.TextFrame.TextRange

# Paragraph 1: "Reminder:" -> "There’s a conference!:"
# (set to a throwaway value first so the engine performs a full run
# replacement instead of merging common characters into a 2nd run)
$para1 = $tr.Paragraphs(1, 1)
$para1.Text = "placeholder"
$para1 = $tr.Paragraphs(1, 1)
$para1.Text = "There’s a conference!:"

# Paragraph 2: "No class Thursday..." -> "Mark needs a count..."
$para2 = $tr.Paragraphs(2, 1)
$para2.Text = "placeholder"
$para2 = $tr.Paragraphs(2, 1)
$para2.Text = "Mark needs a count of how many people plan on attending on Thurs Oct 30. "

# Paragraph 3: "Test pushed to next Tuesday..." is removed entirely.
$para3 = $tr.Paragraphs(3, 1)
$para3.Delete()

Write-Host "Edit complete"
